$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - "Save", copying style (format) from G1 (bold header style)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save column values for rows 2-27
$saveValues = @(0,1,1,1,0,0,0,0,1,1,1,0,1,0,1,1,1,1,0,0,0,1,1,0,1,0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
